# Apply "Update with Correct Forecast output" changes.
#
# Sheet "Forecast Comparison":
#   - Insert a new column B named "Week_Start_Date" holding each week's
#     start date (text), shifting the old B:I columns to C:J.
#   - Normalize the "Week" labels in column A from zero-padded (W01..W09)
#     to non-padded (W1..W9) for the first nine weeks.
#   - Refresh the forecast numbers (old column C/MyForecast values, which
#     become the new column D) to their corrected figures.
#   - Re-type the "is_holiday_week" column (now J) as boolean.
#
# Sheet "Summary":
#   - Update "Total Forecast (4 Weeks)" from 7 to 6.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Forecast Comparison")

# --- Insert the new "Week_Start_Date" column before the existing column B ---
$ws.Range("B1").EntireColumn.Insert()

# Make sure the new column stores its dates as plain text (so values like
# "2025-01-05" are kept literally instead of being converted to date
# serials), then fill in the header + the sixteen weekly start dates.
$ws.Range("B1:B17").NumberFormat = "@"

$ws.Range("B1").Value = "Week_Start_Date"

$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]
}

# --- Normalize the zero-padded week numbers for weeks 1-9 ---
for ($i = 1; $i -le 9; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "W$i"
}

# --- Refresh the corrected forecast figures ---
# Columns (after the insert): A Week, B Week_Start_Date, C ASIN,
# D MyForecast, E Amazon Mean Forecast, F Amazon P70 Forecast,
# G Amazon P80 Forecast, H Amazon P90 Forecast, I Product Title,
# J is_holiday_week
#
# Row => D, E, F, G, H values
$forecastValues = @{
    2  = @(2, 2, 2, 3, 7)
    3  = @(2, 2, 2, 3, 6)
    4  = @(2, 2, 2, 3, 6)
    5  = @(1, 2, 2, 3, 6)
    6  = @(2, 2, 2, 3, 6)
    7  = @(1, 2, 2, 3, 6)
    8  = @(1, 2, 1, 3, 5)
    9  = @(1, 2, 1, 3, 5)
    10 = @(1, 2, 1, 3, 5)
    11 = @(1, 2, 1, 3, 5)
    12 = @(1, 2, 1, 3, 5)
    13 = @(1, 1, 1, 2, 4)
    14 = @(1, 1, 1, 2, 4)
    15 = @(1, 1, 1, 2, 4)
    16 = @(1, 1, 1, 2, 4)
    17 = @(1, 1, 1, 2, 3)
}

foreach ($row in $forecastValues.Keys) {
    $vals = $forecastValues[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
    $ws.Cells.Item($row, 6).Value = $vals[2]
    $ws.Cells.Item($row, 7).Value = $vals[3]
    $ws.Cells.Item($row, 8).Value = $vals[4]
}

# --- Re-type "is_holiday_week" (column J) values as booleans ---
for ($row = 2; $row -le 17; $row++) {
    $ws.Cells.Item($row, 10).Value = $false
}

# --- Update the Summary sheet's 4-week forecast total ---
# (Kept as text, matching the existing "Value" column convention, so
# format the cell as text first to stop it being reinterpreted as a number.)
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B11").NumberFormat = "@"
$summary.Range("B11").Value = "6"
